$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7250
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 4999.75
$ws.Range("I46").Value = 2999.3333
$ws.Range("J46").Value = 6200
$ws.Range("K46").Value = 8997.999899999999
$ws.Range("L46").Value = 18600
$ws.Range("M46").Value = -8878.999899999999
$ws.Range("N46").Value = -18838

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 4999.75
$ws.Range("I60").Value = 2999.3333
$ws.Range("J60").Value = 6200
$ws.Range("K60").Value = 8997.999899999999
$ws.Range("L60").Value = 18600
$ws.Range("M60").Value = -8513.999899999999
$ws.Range("N60").Value = -19568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6439.385
$ws.Range("I86").Value = 6511.4375
$ws.Range("J86").Value = 6324.1
$ws.Range("K86").Value = 6511.4375
$ws.Range("L86").Value = 6324.1
$ws.Range("M86").Value = -5388.4375
$ws.Range("N86").Value = -8570.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 6439.385
$ws.Range("I89").Value = 6511.4375
$ws.Range("J89").Value = 6324.1
$ws.Range("K89").Value = 32557.1875
$ws.Range("L89").Value = 31620.5
$ws.Range("M89").Value = -26941.1875
$ws.Range("N89").Value = -42852.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 628.4286
$ws.Range("J111").Value = 999
$ws.Range("L111").Value = 2997
$ws.Range("N111").Value = -9131

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 3721.0454
$ws.Range("I125").Value = 3715.2
$ws.Range("K125").Value = 33436.8
$ws.Range("M125").Value = -30976.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2738.8
$ws.Range("I127").Value = 3374.25
$ws.Range("K127").Value = 10122.75
$ws.Range("M127").Value = -5162.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1370.4166
$ws.Range("I129").Value = 763.8570999999999
$ws.Range("K129").Value = 2291.5713
$ws.Range("M129").Value = 2708.4287

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 88000
$ws.Range("J136").Value = 88000
$ws.Range("L136").Value = 88000
$ws.Range("N136").Value = -98200

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 16672026
$ws.Range("I137").Value = 22729168
$ws.Range("J137").Value = 14885.75
$ws.Range("K137").Value = 68187504
$ws.Range("L137").Value = 44657.25
$ws.Range("M137").Value = -68184954
$ws.Range("N137").Value = -49757.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4830.185
$ws.Range("I138").Value = 5690.5557
$ws.Range("J138").Value = 4400
$ws.Range("K138").Value = 17071.6671
$ws.Range("L138").Value = 13200
$ws.Range("M138").Value = -11931.6671
$ws.Range("N138").Value = -23480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 928807.9399999999
$ws.Range("J32").Value = 21869.2
$ws.Range("L32").Value = 21869.2
$ws.Range("N32").Value = -22443.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2269.8
$ws.Range("I63").Value = 1750
$ws.Range("J63").Value = 2399.75
$ws.Range("K63").Value = 1750
$ws.Range("L63").Value = 2399.75
$ws.Range("M63").Value = -1064
$ws.Range("N63").Value = -3771.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2269.8
$ws.Range("I66").Value = 1750
$ws.Range("J66").Value = 2399.75
$ws.Range("K66").Value = 8750
$ws.Range("L66").Value = 11998.75
$ws.Range("M66").Value = -5318
$ws.Range("N66").Value = -18862.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 65000.5
$ws.Range("J134").Value = 65000.5
$ws.Range("L134").Value = 65000.5
$ws.Range("N134").Value = -75140.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3553.2856
$ws.Range("I86").Value = 2191
$ws.Range("J86").Value = 4575
$ws.Range("K86").Value = 2191
$ws.Range("L86").Value = 4575
$ws.Range("M86").Value = -1068
$ws.Range("N86").Value = -6821

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3553.2856
$ws.Range("I89").Value = 2191
$ws.Range("J89").Value = 4575
$ws.Range("K89").Value = 10955
$ws.Range("L89").Value = 22875
$ws.Range("M89").Value = -5339
$ws.Range("N89").Value = -34107

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1639.0416
$ws.Range("I107").Value = 1337.1818
$ws.Range("K107").Value = 1337.1818
$ws.Range("M107").Value = 582.8181999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1349.4348
$ws.Range("I22").Value = 836.13336
$ws.Range("J22").Value = 2311.875
$ws.Range("K22").Value = 836.13336
$ws.Range("L22").Value = 2311.875
$ws.Range("M22").Value = -486.13336
$ws.Range("N22").Value = -3011.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4082.5
$ws.Range("I62").Value = 4033.182
$ws.Range("K62").Value = 4033.182
$ws.Range("M62").Value = -3409.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4082.5
$ws.Range("I65").Value = 4033.182
$ws.Range("K65").Value = 20165.91
$ws.Range("M65").Value = -17045.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2081.6667
$ws.Range("I122").Value = 1466.6154
$ws.Range("K122").Value = 4399.8462
$ws.Range("M122").Value = -1949.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2852.6155
$ws.Range("I140").Value = 2398.6086
$ws.Range("K140").Value = 7195.825800000001
$ws.Range("M140").Value = -2015.825800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3444.842
$ws.Range("I122").Value = 3450.2942
$ws.Range("K122").Value = 10350.8826
$ws.Range("M122").Value = -7900.882599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9824.214
$ws.Range("I132").Value = 5467.6055
$ws.Range("K132").Value = 16402.8165
$ws.Range("M132").Value = -13872.8165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5279.6924
$ws.Range("I40").Value = 4545.25
$ws.Range("K40").Value = 4545.25
$ws.Range("M40").Value = -4409.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2300.3462
$ws.Range("J68").Value = 1984.1111
$ws.Range("L68").Value = 1984.1111
$ws.Range("N68").Value = -3482.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2300.3462
$ws.Range("J71").Value = 1984.1111
$ws.Range("L71").Value = 9920.5555
$ws.Range("N71").Value = -17408.5555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17898.572
$ws.Range("I62").Value = 21458
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 21458
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -20834
$ws.Range("N62").Value = -10248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 17898.572
$ws.Range("I65").Value = 21458
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 107290
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -104170
$ws.Range("N65").Value = -51240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 924.6429000000001
$ws.Range("I107").Value = 912.0833
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2736.2499
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -816.2498999999998
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1499.5714
$ws.Range("I126").Value = 1249.5
$ws.Range("K126").Value = 3748.5
$ws.Range("M126").Value = -1278.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3970718
$ws.Range("I132").Value = 4169055.5
$ws.Range("J132").Value = 3971.5
$ws.Range("K132").Value = 12507166.5
$ws.Range("L132").Value = 11914.5
$ws.Range("M132").Value = -12504636.5
$ws.Range("N132").Value = -16974.5
